# Adds 5 new rows (157-161) to the master-reg_center_device_h sheet, continuing
# the existing data pattern: regcntr_id=10002, device_id incrementing from
# 3000176 to 3000180, lang_code="eng", is_active=TRUE, cr_by="superadmin",
# cr_dtimes="now()", eff_dtimes="now()".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 157
$startDevice = 3000176

for ($i = 0; $i -lt 5; $i++) {
    $row = $startRow + $i
    $deviceId = $startDevice + $i

    $ws.Cells.Item($row, 1).Value = 10002
    $ws.Cells.Item($row, 2).Value = $deviceId
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Update the view to reflect the new scroll position / active selection.
$ws.Application.ActiveWindow.ScrollRow = 152
$ws.Range("B157").Select()
